$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.99999998602528151
$ws.Range("A2").Value = 0.99485489673242589
$ws.Range("A3").Value = 0.97558546551599234
$ws.Range("A4").Value = 0.96723191362921768
$ws.Range("A5").Value = 0.95933820489874833
$ws.Range("A6").Value = 0.94299882588221073
$ws.Range("A7").Value = 0.94056228047185853
$ws.Range("A8").Value = 0.93587699582701867
$ws.Range("A9").Value = 0.93264035360537512
$ws.Range("A10").Value = 0.93032657260564799
$ws.Range("A11").Value = 0.9299364599246096
$ws.Range("A12").Value = 0.92945037661783392
$ws.Range("A13").Value = 0.93064618731580473
$ws.Range("A14").Value = 0.92935214111635367
$ws.Range("A15").Value = 0.92953417223145052
$ws.Range("A16").Value = 0.93035873870604313
$ws.Range("A17").Value = 0.92665110552359153
$ws.Range("A18").Value = 0.92554223771288768
$ws.Range("A19").Value = 0.99225846342176727
$ws.Range("A20").Value = 0.98514167717851497
$ws.Range("A21").Value = 0.98374323589508306
$ws.Range("A22").Value = 0.98247873155176513
$ws.Range("A23").Value = 0.96855019440803736
$ws.Range("A24").Value = 0.95552880484406699
$ws.Range("A25").Value = 0.94907174621475554
$ws.Range("A26").Value = 0.94038790680273943
$ws.Range("A27").Value = 0.93699771604779136
$ws.Range("A28").Value = 0.92173059547509306
$ws.Range("A29").Value = 0.91107691704131022
$ws.Range("A30").Value = 0.90701412110978019
$ws.Range("A31").Value = 0.90509950836576736
$ws.Range("A32").Value = 0.90342029631497756
$ws.Range("A33").Value = 0.90290030928171938
